$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks before rewriting data (avoids stale refs)
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-11-11 18:26:27'
$ws.Range("B2").Value = 'ChatGPTを用いた当事業部内チャットツールのシステム開発'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5431738'
$ws.Range("G2").Value = 463
$ws.Range("H2").Value = '🔥GPT,ChatGPT ◆ツール,開発'

# Row 3
$ws.Range("A3").Value = '2025-11-11 18:26:27'
$ws.Range("B3").Value = '【急募】AIシステム構築!FirebaseとOpenAI活用の専門家募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5431299'
$ws.Range("G3").Value = 325
$ws.Range("H3").Value = '🔥AI,Ai'

# Row 4
$ws.Range("A4").Value = '2025-11-11 18:26:27'
$ws.Range("B4").Value = '【急募】Cordova必須!スマホアプリ開発支援メンバー募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5431740'
$ws.Range("G4").Value = 175
$ws.Range("H4").Value = '★スマホアプリ ◆開発 ◇アプリ'

# Row 5
$ws.Range("A5").Value = '2025-11-11 18:26:27'
$ws.Range("B5").Value = 'webアプリの開発'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5431673'
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = '◆開発 ◇アプリ'

# Row 6
$ws.Range("A6").Value = '2025-11-11 18:26:27'
$ws.Range("B6").Value = '【急募】Webアプリ開発エンジニア募集!フルリモート可'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5411585'
$ws.Range("G6").Value = 93
$ws.Range("H6").Value = '◆開発 ◇アプリ'

# Row 7
$ws.Range("A7").Value = '2025-11-11 18:26:27'
$ws.Range("B7").Value = '【急募】知的財産関連システムの開発パートナー募集'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5431547'
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = '◆開発'

# Row 8
$ws.Range("A8").Value = '2025-11-11 18:26:27'
$ws.Range("B8").Value = '進行管理およびチームディレクションを担当'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '~ 5,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = '◇管理'

# Row 9
$ws.Range("A9").Value = '2025-11-11 18:26:27'
$ws.Range("B9").Value = '【急募】Laravel12でFortifyを使った2段階認証システムの制作'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5431508'
$ws.Range("G9").Value = 33

# Row 10
$ws.Range("A10").Value = '2025-11-11 18:26:27'
$ws.Range("B10").Value = '〖リモート可〗Delphiエンジニア募集'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5341051'
$ws.Range("G10").Value = 25

# Row 11
$ws.Range("A11").Value = '2025-11-11 18:26:27'
$ws.Range("B11").Value = '【フルリモート】SESエンジニア募集|スキルに応じて30〜40万円/月|複数案件あり・継続前提'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5417644'
$ws.Range("G11").Value = 25

# Row 12
$ws.Range("A12").Value = '2025-11-11 18:26:27'
$ws.Range("B12").Value = '【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5431322'
$ws.Range("G12").Value = 25

# Row 13
$ws.Range("A13").Value = '2025-11-11 18:26:27'
$ws.Range("B13").Value = '【音楽制作】サイケデリックトランスのトラックを作成してくれる方募集'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5432042'
$ws.Range("G13").Value = 18

# Row 14
$ws.Range("A14").Value = '2025-11-11 18:26:27'
$ws.Range("B14").Value = '初回 Hubspot構築者募集'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5431947'
$ws.Range("G14").Value = 18

# Row 15
$ws.Range("A15").Value = '2025-11-11 18:26:27'
$ws.Range("B15").Value = 'AWS環境からAWS環境ヘの新規構築'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5431069'
$ws.Range("G15").Value = 18

# Row 16
$ws.Range("A16").Value = '2025-11-11 18:26:27'
$ws.Range("B16").Value = '【Stable Diffusion】参考動画に沿って約100プロンプト構築'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5432055'
$ws.Range("G16").Value = 10

# Row 17
$ws.Range("A17").Value = '2025-11-11 18:26:27'
$ws.Range("B17").Value = 'EAの作成'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5431276'
$ws.Range("G17").Value = 10

# Re-create hyperlinks for column F (rows 2-17), reusing the Hyperlink style
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5431738') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5431299') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5431740') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5431673') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5411585') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5431547') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5418064') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5431508') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5341051') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5417644') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5431322') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5432042') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5431947') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5431069') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5432055') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5431276') | Out-Null

# Restore the original Hyperlink cell style (Hyperlinks.Add creates a duplicate style otherwise)
$ws.Range("F2:F17").Style = "Hyperlink"
